{"js": "// Replace the multiplication problems in the practice table with the\n// regenerated set of problems (old -> new text, cell by cell).\nconst replacements = [\n  [\"113\u00d76=\", \"285\u00d77=\"],\n  [\"742\u00d76=\", \"178\u00d77=\"],\n  [\"551\u00d73=\", \"246\u00d77=\"],\n  [\"559\u00d74=\", \"578\u00d72=\"],\n  [\"373\u00d79=\", \"390\u00d73=\"],\n  [\"297\u00d72=\", \"700\u00d78=\"],\n  [\"837\u00d79=\", \"895\u00d73=\"],\n  [\"658\u00d78=\", \"445\u00d79=\"],\n  [\"399\u00d75=\", \"667\u00d76=\"],\n  [\"705\u00d78=\", \"183\u00d79=\"],\n  [\"555\u00d75=\", \"202\u00d77=\"],\n  [\"152\u00d79=\", \"102\u00d79=\"],\n  [\"621\u00d76=\", \"934\u00d74=\"],\n  [\"249\u00d78=\", \"850\u00d79=\"],\n  [\"826\u00d77=\", \"936\u00d73=\"],\n  [\"121\u00d78=\", \"636\u00d79=\"],\n  [\"314\u00d72=\", \"520\u00d74=\"],\n  [\"784\u00d73=\", \"972\u00d74=\"],\n  [\"499\u00d74=\", \"530\u00d79=\"],\n  [\"858\u00d75=\", \"694\u00d77=\"],\n  [\"638\u00d72=\", \"121\u00d72=\"],\n  [\"632\u00d74=\", \"750\u00d77=\"],\n  [\"285\u00d76=\", \"105\u00d77=\"],\n  [\"238\u00d77=\", \"228\u00d72=\"],\n  [\"440\u00d74=\", \"286\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication problems in the practice table with the\n# regenerated set of problems (old -> new text, cell by cell).\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"113\u00d76=\", \"285\u00d77=\"),\n  @(\"742\u00d76=\", \"178\u00d77=\"),\n  @(\"551\u00d73=\", \"246\u00d77=\"),\n  @(\"559\u00d74=\", \"578\u00d72=\"),\n  @(\"373\u00d79=\", \"390\u00d73=\"),\n  @(\"297\u00d72=\", \"700\u00d78=\"),\n  @(\"837\u00d79=\", \"895\u00d73=\"),\n  @(\"658\u00d78=\", \"445\u00d79=\"),\n  @(\"399\u00d75=\", \"667\u00d76=\"),\n  @(\"705\u00d78=\", \"183\u00d79=\"),\n  @(\"555\u00d75=\", \"202\u00d77=\"),\n  @(\"152\u00d79=\", \"102\u00d79=\"),\n  @(\"621\u00d76=\", \"934\u00d74=\"),\n  @(\"249\u00d78=\", \"850\u00d79=\"),\n  @(\"826\u00d77=\", \"936\u00d73=\"),\n  @(\"121\u00d78=\", \"636\u00d79=\"),\n  @(\"314\u00d72=\", \"520\u00d74=\"),\n  @(\"784\u00d73=\", \"972\u00d74=\"),\n  @(\"499\u00d74=\", \"530\u00d79=\"),\n  @(\"858\u00d75=\", \"694\u00d77=\"),\n  @(\"638\u00d72=\", \"121\u00d72=\"),\n  @(\"632\u00d74=\", \"750\u00d77=\"),\n  @(\"285\u00d76=\", \"105\u00d77=\"),\n  @(\"238\u00d77=\", \"228\u00d72=\"),\n  @(\"440\u00d74=\", \"286\u00d76=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Replacement.ClearFormatting()\n  $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
